$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the iccid/imei header columns (B, C) as Text format, matching the
# "Text" number format applied to the legalizador sheet so ICCID/IMEI values
# keep their leading spaces / full digit strings.
$ws.Range("B1:C1").NumberFormat = "@"

# New rows of vendor/legalization data appended below the header.
$rows = @(
    @{ idvendedor = 1010014821; iccid = " 602511579633"; imei = "869466060182485"; min = 3224114582; apellido = "OQUENDO"; nombre = "ALEJANDRO"; cedula = 1007524653; tipodoc = "cc" },
    @{ idvendedor = 1010014821; iccid = " 602306826399"; imei = "868651052191260"; min = 3102587311; apellido = "CATAÑO"; nombre = "ANA"; cedula = 43655411; tipodoc = "cc" },
    @{ idvendedor = 1010014821; iccid = " 602202936480"; imei = "868651052190171"; min = 3102585958; apellido = "MUÑOZ"; nombre = "ANDRES"; cedula = 1007524753; tipodoc = "cc" },
    @{ idvendedor = 1010014821; iccid = " 602308742078"; imei = "355689861577154"; min = 3124372604; apellido = "OLARTE"; nombre = "BLANCA"; cedula = 65680215; tipodoc = "cc" }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.idvendedor

    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row.iccid

    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $row.imei

    $ws.Cells.Item($r, 4).Value = $row.min
    $ws.Cells.Item($r, 5).Value = $row.apellido
    $ws.Cells.Item($r, 6).Value = $row.nombre
    $ws.Cells.Item($r, 7).Value = $row.cedula
    $ws.Cells.Item($r, 8).Value = $row.tipodoc

    $r = $r + 1
}

# Column sizing to match the legalizador template layout (idvendedor column
# auto-fit to its content; iccid/imei columns kept at the default width but
# flagged with the Text number format above).
$ws.Columns("A:A").ColumnWidth = 10.45
$ws.Columns("B:C").ColumnWidth = 8.3

$null = $ws.Range("H5").Select()
